$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = 109714023
$ws.Range("B7").Value = 89790
$ws.Range("C7").Value = "Ovaliderad"
$ws.Range("D7").Value = "NT"
$ws.Range("E7").Value = 6040186
$ws.Range("G7").Value = "Leptoporus mollis"
$ws.Range("H7").Value = "(Pers.:Fr.) Quél."
$ws.Range("I7").Value = ""
$ws.Range("J7").Value = ""
$ws.Range("K7").Value = ""
$ws.Range("N7").Value = ""
$ws.Range("P7").Value = "Dammsjön, Gstr"
$ws.Range("Q7").Value = 575782.5865376759
$ws.Range("R7").Value = 6703744.008187429
$ws.Range("S7").Value = 25
$ws.Range("T7").Value = "Gävleborg"
$ws.Range("U7").Value = "Hofors"
$ws.Range("V7").Value = "Gästrikland"
$ws.Range("W7").Value = "Torsåker"

# Startdatum / Slutdatum look like ISO dates ("2023-06-02"); Excel's COM layer
# auto-parses such literals into date serials when assigned through .Value.
# Source data stores them as plain text, so force text formatting first and
# drop the formatting override afterwards (matches default/no explicit style).
$ws.Range("Y7").NumberFormat = "@"
$ws.Range("Y7").Value = "2023-06-02"
$ws.Range("Y7").ClearFormats()

$ws.Range("Z7").Value = "07:49"

$ws.Range("AA7").NumberFormat = "@"
$ws.Range("AA7").Value = "2023-06-02"
$ws.Range("AA7").ClearFormats()

$ws.Range("AB7").Value = "07:49"
$ws.Range("AC7").Value = "På gran"
$ws.Range("AD7").Value = $false
$ws.Range("AE7").Value = $false
$ws.Range("AF7").Value = ""
$ws.Range("AG7").Value = $false
$ws.Range("AT7").Value = ""
$ws.Range("AW7").Value = "Philipp Weiss"
$ws.Range("AX7").Value = "Philipp Weiss"
$ws.Range("AY7").Value = ""
